$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, copying the style/format of the existing
# header cell G1 (bold, centered, bordered) so it matches the other headers.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the corresponding numeric value in H2 (plain numeric cell, no special style).
$ws.Range("H2").Value = 0
